$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh - GitHub Actions scheduled update

$ws.Range("D2").Value = "43.073.37"
$ws.Range("E2").Value = "  +2.00%  "

$ws.Range("D3").Value = "2.310.72"
$ws.Range("E3").Value = "  +1.91%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.20"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.57"
$ws.Range("E6").Value = "  +5.93%  "

$ws.Range("E7").Value = "  +2.17%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.91"
$ws.Range("E10").Value = "  +8.41%  "

$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.91"
$ws.Range("E13").Value = "  +14.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  +3.76%  "

$ws.Range("D15").Value = "2.687.21"
$ws.Range("E15").Value = "  +2.57%  "

$ws.Range("D16").Value = "2.338.88"
$ws.Range("E16").Value = "  +2.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.816"
$ws.Range("E17").Value = "  +4.25%  "

$ws.Range("D18").Value = "42.990.92"
$ws.Range("E18").Value = "  +1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  +8.60%  "

$ws.Range("E20").Value = "  +3.27%  "

$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +1.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.84"
$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.21"
$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +13.14%  "

$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.81"
$ws.Range("E27").Value = "  +3.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  +1.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.96"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.53"
$ws.Range("E30").Value = "  +2.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.26"
$ws.Range("E31").Value = "  +1.10%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("E33").Value = "  +3.26%  "

$ws.Range("E34").Value = "  +3.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.18"
$ws.Range("E35").Value = "  +3.56%  "

$ws.Range("E36").Value = "  +3.56%  "

$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  +4.58%  "

$ws.Range("E41").Value = "  +1.63%  "

$ws.Range("E42").Value = "  -4.02%  "

$ws.Range("D43").Value = "1.991.56"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0289"
$ws.Range("E44").Value = "  +4.49%  "

$ws.Range("E45").Value = "  +8.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.60"
$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("E47").Value = "  +4.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.32"
$ws.Range("E48").Value = "  +7.42%  "

$ws.Range("D49").Value = "2.530.34"

$ws.Range("E50").Value = "  +3.75%  "

$ws.Range("E51").Value = "  +2.25%  "
